$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.795.55'
$ws.Range("E2").Value = '  +4.93%  '

$ws.Range("D3").Value = '1.889.13'
$ws.Range("E3").Value = '  +3.58%  '

$ws.Range("D4").Value = '''0.9983'
$ws.Range("E4").Value = '  -0.83%  '

$ws.Range("D5").Value = '''338.99'
$ws.Range("E5").Value = '  +2.10%  '

$ws.Range("D6").Value = '''0.9987'
$ws.Range("E6").Value = '  -0.67%  '

$ws.Range("D7").Value = '''0.4730'
$ws.Range("E7").Value = '  +3.35%  '

$ws.Range("D8").Value = '''0.4041'
$ws.Range("E8").Value = '  +6.02%  '

$ws.Range("D9").Value = '''47.61'
$ws.Range("E9").Value = '  +2.67%  '

$ws.Range("D10").Value = '''0.08078'
$ws.Range("E10").Value = '  +2.50%  '

$ws.Range("D11").Value = '''1.012'
$ws.Range("E11").Value = '  +4.56%  '

$ws.Range("D12").Value = '''22.32'
$ws.Range("E12").Value = '  +6.34%  '

$ws.Range("D13").Value = '''6.071'
$ws.Range("E13").Value = '  +3.29%  '

$ws.Range("D14").Value = '1.871.68'
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").Value = '''7.324'
$ws.Range("E15").Value = '  +4.15%  '

$ws.Range("D16").Value = '''90.98'
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").Value = '''0.9990'
$ws.Range("E17").Value = '  -0.78%  '

$ws.Range("D18").Value = '''0.00001048'
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("D19").Value = '''0.06617'
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").Value = '''17.72'
$ws.Range("E20").Value = '  +3.72%  '

$ws.Range("D21").Value = '''0.9997'
$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").Value = '28.797.03'
$ws.Range("E22").Value = '  +4.98%  '

$ws.Range("D23").Value = '''5.521'
$ws.Range("E23").Value = '  +3.55%  '

$ws.Range("D24").Value = '''11.09'
$ws.Range("E24").Value = '  +2.76%  '

$ws.Range("D25").Value = '''2.263'
$ws.Range("E25").Value = '  -1.62%  '

$ws.Range("D26").Value = '2.098.27'
$ws.Range("E26").Value = '  +2.29%  '

$ws.Range("D27").Value = '''160.67'
$ws.Range("E27").Value = '  +2.98%  '

$ws.Range("D28").Value = '''19.91'
$ws.Range("E28").Value = '  +2.97%  '

$ws.Range("D29").Value = '''2.143'
$ws.Range("E29").Value = '  +4.33%  '

$ws.Range("D30").Value = '''5.520'
$ws.Range("E30").Value = '  +4.87%  '

$ws.Range("D31").Value = '''120.15'
$ws.Range("E31").Value = '  +1.69%  '

$ws.Range("D32").Value = '''0.9960'
$ws.Range("E32").Value = '  +5.44%  '

$ws.Range("D33").Value = '''0.09563'
$ws.Range("E33").Value = '  +2.66%  '

$ws.Range("E34").Value = '  +2.26%  '

$ws.Range("D35").Value = '''1.400'
$ws.Range("E35").Value = '  +5.97%  '

$ws.Range("D36").Value = '''5.392'
$ws.Range("E36").Value = '  +2.95%  '

$ws.Range("D37").Value = '''0.06199'
$ws.Range("E37").Value = '  +4.75%  '

$ws.Range("D38").Value = '''0.02281'
$ws.Range("E38").Value = '  +4.37%  '

$ws.Range("D39").Value = '''8.534'
$ws.Range("E39").Value = '  +6.65%  '

$ws.Range("D40").Value = '''1.191'
$ws.Range("E40").Value = '  +2.93%  '

$ws.Range("D41").Value = '''0.5987'
$ws.Range("E41").Value = '  +4.04%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.1895'
$ws.Range("E42").Value = '  +3.57%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '''0.9990'
$ws.Range("E43").Value = '  -0.50%  '

$ws.Range("D44").Value = '''10.42'
$ws.Range("E44").Value = '  +4.13%  '

$ws.Range("D45").Value = '''1.265'
$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").Value = '''0.5611'
$ws.Range("E46").Value = '  +3.16%  '

$ws.Range("D47").Value = '''12.21'
$ws.Range("E47").Value = '  +2.22%  '

$ws.Range("D48").Value = '''1.969'
$ws.Range("E48").Value = '  +5.61%  '

$ws.Range("D49").Value = '''0.07225'
$ws.Range("E49").Value = '  +9.27%  '

$ws.Range("D50").Value = '''2.118'
$ws.Range("E50").Value = '  +14.66%  '

$ws.Range("D51").Value = '''112.72'
$ws.Range("E51").Value = '  +2.13%  '
